# Update the "Tiers" variable sheet: rename header labels and swap
# the two data rows so the header row (nomTiers / clientProspect)
# comes first, followed by the sample row (April / Client).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nomTiers"
$ws.Range("B1").Value = "clientProspect"
$ws.Range("A2").Value = "April"
$ws.Range("B2").Value = "Client"

# Widen column B to fit the new "clientProspect" header text.
$ws.Columns("B").ColumnWidth = 18

# Move/restore the active selection to C10, matching the saved view state.
$null = $ws.Range("C10").Select()
